# Apply cryptos-list price/volume updates (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.408.84'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.878.83'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7175'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07975'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08119'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.15%  '
$ws.Range('D12').Value = '1.887.88'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '95.09'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.240'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7083'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.410'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008441'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').Value = '29.410.73'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '253.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.41%  '
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').Value = '2.135.07'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.692'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1584'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.067'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.28%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.423'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.321'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.226'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05325'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.953'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7582'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.176'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('E37').Value = '  +0.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01895'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').Value = '1.274.77'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.764'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.410'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.23%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.65'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9057'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '111.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -1.66%  '
$ws.Range('D47').Value = '2.031.22'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.809'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.537'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4345'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.25%  '
